$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, matching the source data's storage as
# inline/shared strings (even for values that look numeric, e.g. "250.17" or
# "37.161.46", and percentage strings like "  +2.35%  "). Forcing the
# NumberFormat to "@" (Text) before assigning prevents Excel's normal
# text-to-number/percentage auto-conversion from kicking in.
function Set-CellText($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-CellText "D2" "37.161.46"
Set-CellText "E2" "  +2.35%  "
Set-CellText "D3" "2.086.10"
Set-CellText "E3" "  +3.26%  "
Set-CellText "E4" "  -0.06%  "
Set-CellText "D5" "250.17"
Set-CellText "E5" "  +2.38%  "
Set-CellText "D6" "0.662"
Set-CellText "E6" "  +0.68%  "
Set-CellText "E7" "  -0.04%  "
Set-CellText "D8" "53.90"
Set-CellText "E8" "  +22.13%  "
Set-CellText "D9" "61.70"
Set-CellText "E9" "  +1.86%  "
Set-CellText "E10" "  +4.43%  "
Set-CellText "D11" "0.0743"
Set-CellText "E11" "  +4.46%  "
Set-CellText "E12" "  +7.85%  "
Set-CellText "D13" "15.03"
Set-CellText "E13" "  +6.15%  "
Set-CellText "D14" "2.390.38"
Set-CellText "E14" "  +3.30%  "
Set-CellText "E15" "  +3.80%  "
Set-CellText "D16" "2.084.48"
Set-CellText "E16" "  +3.19%  "
Set-CellText "D17" "5.16"
Set-CellText "E17" "  +6.31%  "
Set-CellText "D18" "37.136.28"
Set-CellText "E18" "  +2.09%  "
Set-CellText "D19" "72.64"
Set-CellText "E19" "  +2.66%  "
Set-CellText "D20" "14.58"
Set-CellText "E20" "  +15.67%  "
Set-CellText "D21" "0.0₃0844"
Set-CellText "E21" "  +4.70%  "
Set-CellText "D22" "240.31"
Set-CellText "E22" "  +2.30%  "
Set-CellText "D23" "5.19"
Set-CellText "E23" "  +6.92%  "
Set-CellText "E24" "  +0.00%  "
Set-CellText "D25" "2.46"
Set-CellText "E25" "  +1.37%  "
Set-CellText "D26" "172.19"
Set-CellText "E26" "  +2.85%  "
Set-CellText "D27" "9.21"
Set-CellText "E27" "  +6.96%  "
Set-CellText "D28" "20.68"
Set-CellText "E28" "  +5.88%  "
Set-CellText "E29" "  +3.45%  "
Set-CellText "E30" "  +2.70%  "
Set-CellText "D31" "22.76"
Set-CellText "E31" "  +5.66%  "
Set-CellText "E32" "  +29.00%  "
Set-CellText "E33" "  +4.83%  "
Set-CellText "E34" "  +7.00%  "
Set-CellText "D35" "0.0902"
Set-CellText "E35" "  +5.16%  "
Set-CellText "E36" "  -0.07%  "
Set-CellText "D37" "4.14"
Set-CellText "E37" "  +4.95%  "
Set-CellText "D38" "1.81"
Set-CellText "E38" "  -3.00%  "
Set-CellText "E39" "  +4.63%  "
Set-CellText "E40" "  +2.42%  "
Set-CellText "D41" "18.15"
Set-CellText "E41" "  +17.12%  "
Set-CellText "E42" "  +5.75%  "
Set-CellText "E43" "  +5.74%  "
Set-CellText "B44" "FTXToken"
Set-CellText "C44" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-CellText "D44" "4.49"
Set-CellText "E44" "  +131.06%  "
Set-CellText "B45" "Aave"
Set-CellText "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText "D45" "98.25"
Set-CellText "E45" "  +3.31%  "
Set-CellText "B46" "Cronos"
Set-CellText "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText "D46" "0.0943"
Set-CellText "E46" "  +15.84%  "
Set-CellText "E47" "  -0.35%  "
Set-CellText "D48" "1.320.59"
Set-CellText "E48" "  +1.49%  "
Set-CellText "E49" "  +5.32%  "
Set-CellText "D50" "6.97"
Set-CellText "E50" "  +15.20%  "
Set-CellText "D51" "2.32"
Set-CellText "E51" "  +6.33%  "
